$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.02131377
$ws.Range("D3").Value = 0.04708746

Write-Host "ChartObjects count: $($ws.ChartObjects().Count)"
$chart = $ws.ChartObjects(1).Chart
Write-Host "Chart: $chart"
$series = $chart.SeriesCollection(1)
Write-Host "Series: $series"
$series.Values = $ws.Range("D2:D8")
Write-Host "Done"
